$wb = $excel.ActiveWorkbook

# --- DatosCuenta (sheet1): update applicant name / doc number ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "PMSmokeT"
$wsCuenta.Range("B2").Value = "ApellidoPMSmokeT"
$wsCuenta.Range("C2").Value = 27100106

# --- DatosHogar (sheet2): bump household number ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 628

# --- DatosMotor (sheet3): bump vehicle identifiers ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA009"
$wsMotor.Range("B2").Value = "ABC12SSMA009"
$wsMotor.Range("C2").Value = "ZAZ123SSMA009"

# --- DatosAP (sheet4): bump AP number ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200108

# --- Update the saved view state: DatosCuenta is no longer the active tab,
#     its selection moves to C3; DatosAP becomes active with selection E9 ---
$wsCuenta.Activate()
$wsCuenta.Range("C3").Select()

$wsAP.Activate()
$wsAP.Range("E9").Select()
